$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header in H1, matching the header style used by G1 (bold/centered/bordered)
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# Fill H2:H52 with 1 if G>9 else 0 (era/save data)
for ($row = 2; $row -le 52; $row++) {
    $g = $ws.Cells.Item($row, 7).Value2
    if ($g -gt 9) {
        $save = 1
    } else {
        $save = 0
    }
    $ws.Cells.Item($row, 8).Value = $save
}
